# Adds a new "Final Project" assignment block (Final Project / Grade / Comments)
# after the existing "12 CPP" block, and normalizes the border styling of the
# previously-added "10 CPP" / "11 CPP" / "12 CPP" blocks so they match the rest
# of the sheet (thin-bordered cells instead of borderless ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$lastRow = 18

# Column numbers (1-based) for the relevant blocks.
$colAK = 37   # existing separator template (header row style "2")
$colAL = 38   # existing data template column 1 (Assignment header / data)
$colAM = 39   # existing data template column 2 (Grade header / data)
$colAN = 40   # existing data template column 3 (Comments header / data)

$colAO = 41   # new separator before "10 CPP" block
$colAP = 42   # "10 CPP" assignment column (already existed, style fix only)
$colAQ = 43   # "10 CPP" grade column (already existed, style fix only)
$colAR = 44   # "10 CPP" comments column (already existed, style fix only)

$colAS = 45   # new separator before "11 CPP" block
$colAT = 46   # "11 CPP" assignment column (already existed, style fix only)
$colAU = 47   # "11 CPP" grade column (already existed, style fix only)
$colAV = 48   # "11 CPP" comments column (already existed, style fix only)

$colAW = 49   # new separator before "12 CPP" block
$colAX = 50   # "12 CPP" assignment column (already existed, style fix only)
$colAY = 51   # "12 CPP" grade column (already existed, style fix only)
$colAZ = 52   # "12 CPP" comments column (brand new, blank)

$colBA = 53   # new separator before "Final Project" block
$colBB = 54   # "Final Project" assignment column (new)
$colBC = 55   # "Final Project" grade column (new)
$colBD = 56   # "Final Project" comments column (new)

$newSeparatorCols = @($colAO, $colAS, $colAW, $colBA)
$restyleDataCols = @($colAP, $colAQ, $colAR, $colAT, $colAU, $colAV, $colAX, $colAY)

# --- Header row (row 1): build the 4th/new block headers by cloning the
# formatting of the last existing block (AK:AN) and writing the standard
# Assignment / Grade / Comments labels.
foreach ($sepCol in $newSeparatorCols) {
    $ws.Cells.Item(1, $colAK).Copy()
    $ws.Cells.Item(1, $sepCol).PasteSpecial($xlPasteFormats)
}

$headerTriples = @(
    @($colAP, $colAQ, $colAR),
    @($colAT, $colAU, $colAV),
    @($colAX, $colAY, $colAZ),
    @($colBB, $colBC, $colBD)
)

foreach ($triple in $headerTriples) {
    $ws.Cells.Item(1, $colAL).Copy()
    $ws.Cells.Item(1, $triple[0]).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item(1, $triple[0]).Value = "Assignment"

    $ws.Cells.Item(1, $colAM).Copy()
    $ws.Cells.Item(1, $triple[1]).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item(1, $triple[1]).Value = "Grade"

    $ws.Cells.Item(1, $colAN).Copy()
    $ws.Cells.Item(1, $triple[2]).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item(1, $triple[2]).Value = "Comments"
}

# --- Row 2: fix up the border styling on the pre-existing "10/11/12 CPP"
# value cells (they were borderless; make them match the bordered style used
# everywhere else on the row), add the new separators, the new blank "12 CPP"
# comments cell, and the new "Final Project" data.
foreach ($col in $restyleDataCols) {
    $ws.Cells.Item(2, $colAL).Copy()
    $ws.Cells.Item(2, $col).PasteSpecial($xlPasteFormats)
}

foreach ($sepCol in $newSeparatorCols) {
    $ws.Cells.Item(1, $colAK).Copy()
    $ws.Cells.Item(2, $sepCol).PasteSpecial($xlPasteFormats)
}

$ws.Cells.Item(2, $colAL).Copy()
$ws.Cells.Item(2, $colAZ).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(2, $colAL).Copy()
$ws.Cells.Item(2, $colBB).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(2, $colBB).Value = "Final Project"

$ws.Cells.Item(2, $colAM).Copy()
$ws.Cells.Item(2, $colBC).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(2, $colBC).Value = 90

$ws.Cells.Item(2, $colAN).Copy()
$ws.Cells.Item(2, $colBD).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(2, $colBD).Value = "Very Good!"

# --- Rows 3-18: extend the blank, bordered grid out to the new columns.
for ($r = 3; $r -le $lastRow; $r++) {
    foreach ($sepCol in $newSeparatorCols) {
        $ws.Cells.Item(1, $colAK).Copy()
        $ws.Cells.Item($r, $sepCol).PasteSpecial($xlPasteFormats)
    }
    foreach ($triple in $headerTriples) {
        $ws.Cells.Item($r, $colAL).Copy()
        $ws.Cells.Item($r, $triple[0]).PasteSpecial($xlPasteFormats)

        $ws.Cells.Item($r, $colAM).Copy()
        $ws.Cells.Item($r, $triple[1]).PasteSpecial($xlPasteFormats)

        $ws.Cells.Item($r, $colAN).Copy()
        $ws.Cells.Item($r, $triple[2]).PasteSpecial($xlPasteFormats)
    }
}

# --- View state: select the new "Final Project" comments cell range and
# scroll the window so the new block is visible, matching the author's
# saved view (selection = BB2:BD2, top-left visible column = AH).
$ws.Range("BB2:BD2").Select() | Out-Null
$win = $excel.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 34
